$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextPercent($addr, $text) {
    # Percent-looking strings ("74.40%") get auto-coerced into numbers by a
    # plain Range.Value assignment. Route them through a quoted formula and
    # then "Paste Values" over themselves so the stored cell ends up as a
    # plain shared-string text value (matching the original authoring),
    # while keeping the cell's existing style untouched.
    $ws.Range($addr).Formula = '="' + $text + '"'
    $ws.Range($addr).Copy()
    $ws.Range($addr).PasteSpecial(-4163)
}

# --- Column F now carries the CPU info that used to live only in column G ---
$ws.Range("F1").Value = "CPU: Intel(R) Core(TM) i7-8700 CPU @ 3.20GHz, CPU family: x86_64, GPU: Undefined, RAM size: 65702408 kB, OS family: Linux, OS version: Linux-5.15.0-84-generic-x86_64-with-glibc2.29, Python version: 3.8.10"
$ws.Range("F2").Value = "TVM"
$ws.Range("F3").Value = "CPU"
$ws.Range("F4").Value = "FP32"

# --- densenet-121-tf block: refreshed accuracy numbers ---
$ws.Range("A5").Value = "classification"
$ws.Range("B5").Value = "densenet-121-tf"
$ws.Range("C5").Value = "TensorFlow"
$ws.Range("D5").Value = "ImageNet"
$ws.Range("E5").Value = "accuracy@top1"
Set-TextPercent "F5" "74.40%"
$ws.Range("E6").Value = "accuracy@top5"
Set-TextPercent "F6" "92.14%"

# --- googlenet-v1 block moves up into what used to be the efficientnet-b0 slot ---
$ws.Range("B7").Value = "googlenet-v1"
$ws.Range("C7").Value = "Caffe"
$ws.Range("D7").Value = "ImageNet"
$ws.Range("E7").Value = "accuracy@top1"
Set-TextPercent "F7" "68.92%"
$ws.Range("E8").Value = "accuracy@top5"
Set-TextPercent "F8" "89.14%"

# --- googlenet-v4-tf block moves up one slot ---
$ws.Range("B9").Value = "googlenet-v4-tf"
$ws.Range("C9").Value = "TensorFlow"
$ws.Range("D9").Value = "ImageNet"
$ws.Range("E9").Value = "accuracy@top1"
Set-TextPercent "F9" "80.21%"
$ws.Range("E10").Value = "accuracy@top5"
Set-TextPercent "F10" "95.19%"

# --- resnet50-pytorch block moves up one slot ---
$ws.Range("B11").Value = "resnet50-pytorch"
$ws.Range("C11").Value = "PyTorch"
$ws.Range("D11").Value = "ImageNet"
$ws.Range("E11").Value = "accuracy@top1"
Set-TextPercent "F11" "76.15%"
$ws.Range("E12").Value = "accuracy@top5"
Set-TextPercent "F12" "92.87%"

# --- squeezenet1.1 block moves up one slot, refreshed accuracy numbers ---
$ws.Range("B13").Value = "squeezenet1.1"
$ws.Range("C13").Value = "Caffe"
$ws.Range("D13").Value = "ImageNet"
$ws.Range("E13").Value = "accuracy@top1"
Set-TextPercent "F13" "58.38%"
$ws.Range("E14").Value = "accuracy@top5"
Set-TextPercent "F14" "81.01%"

# --- efficientnet-b0 block moves down to the last slot, refreshed accuracy numbers ---
$ws.Range("B15").Value = "efficientnet-b0"
$ws.Range("C15").Value = "TensorFlow"
$ws.Range("D15").Value = "ImageNet"
$ws.Range("E15").Value = "accuracy@top1"
Set-TextPercent "F15" "75.69%"
$ws.Range("E16").Value = "accuracy@top5"
Set-TextPercent "F16" "92.76%"

# --- Column G (the stale duplicate CPU column) is removed entirely ---
$ws.Columns("G").Delete()

# Deleting the column leaves its conditional-formatting rules behind
# (they still reference the now-empty G range) - drop those too.
$ws.Range("G1:G16").FormatConditions.Delete()
